# Swap the presentation's applied theme palette back to the default
# "Office Theme" colours (the deck currently carries the "Integral" /
# "Red Violet" palette on its one-and-only Slide Master).
#
# PowerPoint stores a Theme's 12 colour slots as a ThemeColorScheme on
# the Master; RGB values are plain `long`s in 0x00BBGGRR form (the
# same packing `RGB()` uses), so "RRGGBB" -> R | (G<<8) | (B<<16).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# msoThemeColorDark1 .. msoThemeColorFollowedHyperlink (indices 1-12)
$colorScheme.Item(1).RGB  = 0         # Dark 1    - 000000
$colorScheme.Item(2).RGB  = 16777215  # Light 1   - FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # Dark 2    - 44546A
$colorScheme.Item(4).RGB  = 15132391  # Light 2   - E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # Accent 1  - 5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # Accent 2  - ED7D31
$colorScheme.Item(7).RGB  = 10855845  # Accent 3  - A5A5A5
$colorScheme.Item(8).RGB  = 49407     # Accent 4  - FFC000
$colorScheme.Item(9).RGB  = 12874308  # Accent 5  - 4472C4
$colorScheme.Item(10).RGB = 4697456   # Accent 6  - 70AD47
$colorScheme.Item(11).RGB = 12673797  # Hyperlink - 0563C1
$colorScheme.Item(12).RGB = 7491477   # Followed Hyperlink - 954F72
